$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text runs via Characters, preserves per-run formatting) ---
$ws.Range("A8").Characters(21, 2).Text = "49"
$ws.Range("C9").Characters(27, 10).Text = "12/2/2024"
$ws.Range("C9").Characters(47, 9).Text = "12/8/2024"

# --- Column E width ---
$ws.Columns("E").ColumnWidth = 7.433768

# --- Style template cells (untouched by this edit, used as PasteSpecial format sources) ---
$styleTemplate13 = $ws.Range("A14")
$styleTemplate14 = $ws.Range("I14")
$styleTemplate15 = $ws.Range("K14")

# --- Data table cell updates ---
$ws.Range("G15").Value = "'0"
$styleTemplate13.Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = "'***.*"
$styleTemplate13.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("D16").Value = "'0"
$styleTemplate13.Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "'***.*"
$styleTemplate13.Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 80
$ws.Range("I16").Value = 112
$ws.Range("K16").Value = -15.789473684210
$ws.Range("L16").Value = -36.723163841807
$ws.Range("M16").Value = -54.471544715447
$ws.Range("N16").Value = -88.357588357588
$ws.Range("C17").Value = "'0"
$styleTemplate13.Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -63.636363636363
$ws.Range("I17").Value = 259
$ws.Range("J17").Value = 235
$ws.Range("K17").Value = 10.212765957446
$ws.Range("L17").Value = 30.150753768844
$ws.Range("M17").Value = 175.531914893617
$ws.Range("N17").Value = 14.096916299559
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 213
$ws.Range("J18").Value = 251
$ws.Range("K18").Value = -15.139442231075
$ws.Range("L18").Value = -4.910714285714
$ws.Range("M18").Value = -20.224719101123
$ws.Range("N18").Value = -84.198813056379
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -30.612244897959
$ws.Range("I19").Value = 509
$ws.Range("J19").Value = 592
$ws.Range("K19").Value = -14.020270270270
$ws.Range("L19").Value = -18.035426731078
$ws.Range("M19").Value = 19.203747072599
$ws.Range("N19").Value = -10.858143607705
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 123.076923076923
$ws.Range("I20").Value = 363
$ws.Range("J20").Value = 320
$ws.Range("K20").Value = 13.4375
$ws.Range("L20").Value = 42.352941176470
$ws.Range("M20").Value = 72.037914691943
$ws.Range("N20").Value = -92.191869219186
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -8
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -8.035714285714
$ws.Range("I21").Value = 1480
$ws.Range("J21").Value = 1546
$ws.Range("K21").Value = -4.269081500646
$ws.Range("L21").Value = -1.135604542418
$ws.Range("M21").Value = 17.274167987321
$ws.Range("N21").Value = -81.006160164271
$ws.Range("C22").Value = "'0"
$styleTemplate13.Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("L22").Value = -2.777777777777
$ws.Range("C23").Value = "'0"
$styleTemplate13.Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").Value = "'0"
$styleTemplate13.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "'***.*"
$styleTemplate13.Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 50
$ws.Range("L23").Value = -7.272727272727
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -5.405405405405
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = -6.930693069306
$ws.Range("I24").Value = 1130
$ws.Range("J24").Value = 1342
$ws.Range("K24").Value = -15.797317436661
$ws.Range("L24").Value = -20.478536242083
$ws.Range("M24").Value = 21.505376344086
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -57.142857142857
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = -42.5
$ws.Range("I25").Value = 464
$ws.Range("J25").Value = 499
$ws.Range("K25").Value = -7.014028056112
$ws.Range("L25").Value = 8.411214953271
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 57.142857142857
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 38
$ws.Range("I26").Value = 493
$ws.Range("J26").Value = 453
$ws.Range("K26").Value = 8.830022075055
$ws.Range("L26").Value = 17.102137767220
$ws.Range("M26").Value = 29.396325459317
$ws.Range("G27").Value = "'0"
$styleTemplate13.Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H27").Value = "'***.*"
$styleTemplate13.Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("L27").Value = 3.333333333333
$ws.Range("C28").Value = 4
$styleTemplate14.Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$styleTemplate14.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = 300
$styleTemplate15.Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 47
$ws.Range("K28").Value = -12.765957446808
$ws.Range("L28").Value = -25.454545454545
$ws.Range("D31").Value = 1
$styleTemplate14.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$styleTemplate15.Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("J31").Value = 26
$ws.Range("K31").Value = -73.076923076923
$ws.Range("L33").Value = -16.666666666666
